$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new "Save" header in H1, copying the formatting used by the other
# header cells (e.g. G1 - bold, bordered, centered).
$ws.Range("H1").Value = "Save"
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats

# Fill the Save column with 1 for each data row
$ws.Range("H2").Value = 1
$ws.Range("H3").Value = 1
$ws.Range("H4").Value = 1
